$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2:G7').Value = 'hiragana gojuon seion a_series'
$ws.Range('G8:G11').Value = 'hiragana gojuon seion k_series'
$ws.Range('G12:G16').Value = 'hiragana gojuon seion s_series'
$ws.Range('G17:G21').Value = 'hiragana gojuon seion t_series'
$ws.Range('G22:G26').Value = 'hiragana gojuon seion n_series'
$ws.Range('G27:G31').Value = 'hiragana gojuon seion h_series'
$ws.Range('G32:G36').Value = 'hiragana gojuon seion m_series'
$ws.Range('G37:G39').Value = 'hiragana gojuon seion y_series'
$ws.Range('G40:G44').Value = 'hiragana gojuon seion r_series'
$ws.Range('G45:G47').Value = 'hiragana gojuon seion w_series'
$ws.Range('G48:G52').Value = 'hiragana dakuon k_series g_series'
$ws.Range('G53:G57').Value = 'hiragana dakuon s_series z_series'
$ws.Range('G58:G62').Value = 'hiragana dakuon t_series d_series'
$ws.Range('G63:G67').Value = 'hiragana dakuon h_series b_series'
$ws.Range('G68:G72').Value = 'hiragana dakuon handakuon h_series p_series'
$ws.Range('G73:G75').Value = 'hiragana yoon k_series ky_series'
$ws.Range('G76:G78').Value = 'hiragana yoon s_series sh_series'
$ws.Range('G79:G81').Value = 'hiragana yoon t_series ch_series'
$ws.Range('G82:G84').Value = 'hiragana yoon n_series ny_series'
$ws.Range('G85:G87').Value = 'hiragana yoon h_series hy_series'
$ws.Range('G88:G90').Value = 'hiragana yoon m_series my_series'
$ws.Range('G91:G93').Value = 'hiragana yoon r_series ry_series'
$ws.Range('G94:G96').Value = 'hiragana yoon_dakuon k_series ky_series g_series gy_series'
$ws.Range('G97:G99').Value = 'hiragana yoon_dakuon s_series sh_series j_series'
$ws.Range('G100:G102').Value = 'hiragana yoon_dakuon t_series ch_series j_series'
$ws.Range('G103:G105').Value = 'hiragana yoon_dakuon h_series hy_series b_series by_series'
$ws.Range('G106:G108').Value = 'hiragana yoon_dakuon yoon_handakuon h_series hy_series p_series py_series'
$ws.Range('G109:G113').Value = 'katakana gojuon seion a_series'
$ws.Range('G114:G118').Value = 'katakana gojuon seion k_series'
$ws.Range('G119:G123').Value = 'katakana gojuon seion s_series'
$ws.Range('G124:G128').Value = 'katakana gojuon seion t_series'
$ws.Range('G129:G133').Value = 'katakana gojuon seion n_series'
$ws.Range('G134:G138').Value = 'katakana gojuon seion h_series'
$ws.Range('G139:G143').Value = 'katakana gojuon seion m_series'
$ws.Range('G144:G146').Value = 'katakana gojuon seion y_series'
$ws.Range('G147:G151').Value = 'katakana gojuon seion r_series'
$ws.Range('G152:G154').Value = 'katakana gojuon seion w_series'
$ws.Range('G155:G159').Value = 'katakana dakuon k_series g_series'
$ws.Range('G160:G164').Value = 'katakana dakuon s_series z_series'
$ws.Range('G165:G169').Value = 'katakana dakuon t_series d_series'
$ws.Range('G170:G174').Value = 'katakana dakuon h_series b_series'
$ws.Range('G175:G179').Value = 'katakana dakuon handakuon h_series p_series'
$ws.Range('G180:G182').Value = 'katakana yoon k_series ky_series'
$ws.Range('G183:G185').Value = 'katakana yoon s_series sh_series'
$ws.Range('G186:G188').Value = 'katakana yoon t_series ch_series'
$ws.Range('G189:G191').Value = 'katakana yoon n_series ny_series'
$ws.Range('G192:G194').Value = 'katakana yoon h_series hy_series'
$ws.Range('G195:G197').Value = 'katakana yoon m_series my_series'
$ws.Range('G198:G200').Value = 'katakana yoon r_series ry_series'
$ws.Range('G201:G203').Value = 'katakana yoon_dakuon k_series ky_series g_series gy_series'
$ws.Range('G204:G206').Value = 'katakana yoon_dakuon s_series sh_series j_series'
$ws.Range('G207:G209').Value = 'katakana yoon_dakuon t_series ch_series j_series'
$ws.Range('G210:G212').Value = 'katakana yoon_dakuon h_series hy_series b_series by_series'
$ws.Range('G213:G215').Value = 'katakana yoon_dakuon yoon_handakuon h_series hy_series p_series py_series'
$ws.Range('G216:G256').Value = 'katakana special foreign'

$ws.Columns.Item(7).ColumnWidth = 72.5703125

